$wb = $excel.ActiveWorkbook

# La columna "Tipo" guardaba un codigo numerico para el tipo de relacion.
# Ahora la interfaz grafica trabaja con el texto de la relacion directamente,
# asi que traducimos los codigos existentes a su etiqueta correspondiente.
$tipoPorCodigo = @{
    1 = "Compañero"
    2 = "Conocido"
    3 = "Amigo Personal"
}

# --- Hoja 1 ---
$ws1 = $wb.Worksheets.Item("Hoja 1")
for ($row = 2; $row -le 8; $row++) {
    $codigo = $ws1.Cells.Item($row, 3).Value2
    if ($tipoPorCodigo.ContainsKey([int]$codigo)) {
        $ws1.Cells.Item($row, 3).Value = $tipoPorCodigo[[int]$codigo]
    }
}

# --- Hoja 2 ---
$ws2 = $wb.Worksheets.Item("Hoja 2")

# Correccion de los amigos cargados para Jose y Agustin.
$ws2.Range("B3").Value = "Sebastian"
$ws2.Range("B4").Value = "Juan"

for ($row = 2; $row -le 5; $row++) {
    $codigo = $ws2.Cells.Item($row, 3).Value2
    if ($tipoPorCodigo.ContainsKey([int]$codigo)) {
        $ws2.Cells.Item($row, 3).Value = $tipoPorCodigo[[int]$codigo]
    }
}
